# Generate Report for Archive
#
# The report generator re-ran and the data rows for the two source files
#   1b523678-e84e-4ff7-915b-c4eeb69f95da.md
#   80b10ede-d197-49ee-99de-e8f902fade3a.md
# traded places (row 5 <-> row 6) on every sheet, and 80b10ede's status flipped
# from "Ready for handoff" to "In Translation" (it is now back in translation).
# Hyperlink relationships (their URL targets / r:id) stay anchored to the cell
# they were already attached to - only the visible display text needs to track
# the new cell content.

$wb = $excel.ActiveWorkbook

function Set-CellText {
    param($ws, $ref, $text)
    $ws.Range($ref).Value2 = $text
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $ws.Range($ref).Address()) {
            $h.TextToDisplay = $text
        }
    }
}

function Swap-Row {
    param($ws, $cols, $row1, $row2)
    foreach ($col in $cols) {
        $ref1 = "$col$row1"
        $ref2 = "$col$row2"
        $v1 = $ws.Range($ref1).Value2
        $v2 = $ws.Range($ref2).Value2
        Set-CellText $ws $ref1 $v2
        Set-CellText $ws $ref2 $v1
    }
}

# --- Overview sheet: columns A (file), B (zh-cn status), C (de-de status), D (date)
$wsOverview = $wb.Worksheets.Item("Overview")
Swap-Row $wsOverview @("A","B","C","D") 5 6
# 80b10ede is now on row 5; it went back into translation
Set-CellText $wsOverview "B5" "In Translation"
Set-CellText $wsOverview "C5" "In Translation"

# --- zh-cn sheet: columns A (file), B (ext), C (status), D (handoff file), E (handoff datetime)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Swap-Row $wsZhCn @("A","B","C","D","E") 5 6
Set-CellText $wsZhCn "C5" "In Translation"

# --- de-de sheet: columns A (file), B (ext), C (status), D (handoff file), E (handoff datetime)
$wsDeDe = $wb.Worksheets.Item("de-de")
Swap-Row $wsDeDe @("A","B","C","D","E") 5 6
Set-CellText $wsDeDe "C5" "In Translation"
